$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.955.95"
$ws.Range("E2").Value = "  +6.84%  "
$ws.Range("D3").Value = "3.117.36"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("E4").Value = "  +0.06%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "586.83"
$r.ClearFormats()
$ws.Range("E5").Value = "  +4.72%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "144.18"
$r.ClearFormats()
$ws.Range("E6").Value = "  +5.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.106.29"
$ws.Range("E8").Value = "  +4.04%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.534"
$r.ClearFormats()
$ws.Range("E9").Value = "  +2.34%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.144"
$r.ClearFormats()
$ws.Range("E10").Value = "  +9.31%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "5.74"
$r.ClearFormats()
$ws.Range("E11").Value = "  +11.09%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.472"
$r.ClearFormats()
$ws.Range("E12").Value = "  +3.39%  "
$ws.Range("E13").Value = "  +6.52%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "35.70"
$r.ClearFormats()
$ws.Range("E14").Value = "  +6.28%  "
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "3.632.37"
$ws.Range("E16").Value = "  +4.23%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "7.33"
$r.ClearFormats()
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "3.116.15"
$ws.Range("E18").Value = "  +4.24%  "
$ws.Range("D19").Value = "62.849.74"
$ws.Range("E19").Value = "  +6.62%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "455.60"
$r.ClearFormats()
$ws.Range("E20").Value = "  +6.62%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "14.15"
$r.ClearFormats()
$ws.Range("E21").Value = "  +3.40%  "
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("E23").Value = "  +6.11%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "13.80"
$r.ClearFormats()
$ws.Range("E24").Value = "  +3.98%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "82.54"
$r.ClearFormats()
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +4.85%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "2.70"
$r.ClearFormats()
$ws.Range("E28").Value = "  +6.46%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "8.29"
$r.ClearFormats()
$ws.Range("E29").Value = "  +6.43%  "
$ws.Range("E30").Value = "  +0.18%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "6.87"
$r.ClearFormats()
$ws.Range("E31").Value = "  +13.67%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "0.113"
$r.ClearFormats()
$ws.Range("E32").Value = "  +13.78%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "27.20"
$r.ClearFormats()
$ws.Range("E33").Value = "  +5.91%  "
$ws.Range("E34").Value = "  +5.61%  "
$ws.Range("D35").Value = "0.0₃0810"
$ws.Range("E35").Value = "  +6.79%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "6.12"
$r.ClearFormats()
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("E37").Value = "  +6.95%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "3.07"
$r.ClearFormats()
$ws.Range("E38").Value = "  +12.74%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "50.72"
$r.ClearFormats()
$ws.Range("E39").Value = "  +4.00%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "8.85"
$r.ClearFormats()
$ws.Range("E40").Value = "  +1.96%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "429.26"
$r.ClearFormats()
$ws.Range("E41").Value = "  +7.88%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.937.84"
$ws.Range("E42").Value = "  +6.86%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.0375"
$r.ClearFormats()
$ws.Range("E43").Value = "  +6.71%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.282"
$r.ClearFormats()
$ws.Range("E44").Value = "  +12.39%  "
$ws.Range("E45").Value = "  +2.74%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "2.18"
$r.ClearFormats()
$ws.Range("E46").Value = "  +8.90%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "35.67"
$r.ClearFormats()
$ws.Range("E47").Value = "  +0.23%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "124.55"
$r.ClearFormats()
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("E50").Value = "  +1.77%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "24.74"
$r.ClearFormats()
$ws.Range("E51").Value = "  +6.02%  "